$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Row 91: replace the old shared formula with a direct reference to the new row 92
$ws.Range("D91").Formula = "=D92"

# New row 92 - LZ E5 and CBLS titrations 20220729
$ws.Range("A92").Value = 20220729
$ws.Range("B92").Value = 2226.06306
$ws.Range("C92").Value = 2224.4699999999998
$ws.Range("D92").Formula = "=100*(B92-C92)/C92"
$ws.Range("E92").Value = 180
$ws.Range("F92").Value = "CRM OPENED 20220722 LHZ"

# Update the view to reflect scrolling to the newly added row
$ws.Range("F92").Select()
$excel.ActiveWindow.ScrollRow = 76
